$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster FAPs | Rbp4 | Stra6 | Target cluster -> ECs (new)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7013473333333332
$ws.Range("H2").Value = 2.104042
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.026907
$ws.Range("N2").Value = 0.080721
$ws.Range("O2").Value = 0.02822345482694817
$ws.Range("P2").Value = 0.02822345482694817
$ws.Range("Q2").Value = 0.018871152698
$ws.Range("R2").Value = 0.169840374282
$ws.Range("S2").Value = 0.02822345482694817
$ws.Range("T2").Value = 0.02822345482694817

# Row 3: Target cluster -> FAPs (existing string, was sCs)
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7013473333333332
$ws.Range("H3").Value = 2.104042
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6065346666666667
$ws.Range("N3").Value = 1.819604
$ws.Range("O3").Value = 0.6362100481527012
$ws.Range("P3").Value = 0.6362100481527012
$ws.Range("Q3").Value = 0.4253914710408888
$ws.Range("R3").Value = 3.828523239367999
$ws.Range("S3").Value = 0.6362100481527012
$ws.Range("T3").Value = 0.6362100481527012

# New row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rbp4"
$ws.Range("C4").Value = "Stra6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7013473333333332
$ws.Range("H4").Value = 2.104042
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3199143333333334
$ws.Range("N4").Value = 0.959743
$ws.Range("O4").Value = 0.3355664970203506
$ws.Range("P4").Value = 0.3355664970203506
$ws.Range("Q4").Value = 0.2243710645784444
$ws.Range("R4").Value = 2.019339581206
$ws.Range("S4").Value = 0.3355664970203506
$ws.Range("T4").Value = 0.3355664970203506
